# INVENTARIO.xlsx
#
# Feature: when an item is sold, show how many are left in inventory.
#
# Row 1 ("Coca") had 15 units recorded as sold; after the sale the item's
# "sold" counter resets to reflect that none remain to be (re)sold right
# now, i.e. B1 becomes "0" (matching the other rows of column B, which are
# stored as text, not numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = $ws.Range("B1")

# Writing the literal string "0" straight into B1 would be auto-detected
# as a number (losing the text type that this column already uses for its
# values). To keep it stored as text - same as the rest of column B - we
# stage the text value in a scratch cell, force it to Text format there,
# copy the resulting (already-typed) value onto B1, then discard the
# scratch cell so it leaves no trace in the sheet.
$scratch = $ws.Range("D1")
$scratch.NumberFormat = "@"
$scratch.Value = "0"

$scratch.Copy()
$target.PasteSpecial(-4163)  # xlPasteValues - keeps B1's existing formatting/style

# Remove the helper cell completely (shifting cells left) so the sheet's
# used range / dimensions are left exactly as they were.
$scratch.Delete(-4159)  # xlShiftToLeft
